# StateFuncResources.xlsx - SLG building config update
# Populate column B ("Atlas_ResID" header) with the resource-icon / UI
# setting that belongs to every state-function row in column A.
#   Row 1 (header)         -> B1 already "Atlas_ResID" - leave as is
#   Row 2 (EFT_INFO)       -> "msg_icon"
#   Rows 3-15 (all others) -> "Ssetting"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "msg_icon"

for ($r = 3; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = "Ssetting"
}

# Match the author's last-saved selection.
$ws.Range("E14").Select()
